$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B5 holds the Arabic name for gender code "MLE" (row 5). It previously
# (incorrectly) reused the Arabic word for "Female" (أنثى); fix it to the
# Arabic word for "Male" (الذكر) - a new shared-string entry.
$ws.Range("B5").Value = "الذكر"

# Column B is resized to fit its (now mixed-script) contents.
$ws.Columns("B:B").ColumnWidth = 6.53

# The sheet was left with D16 as the active / selected cell.
$ws.Range("D16").Select()

# Page setup (paper size / orientation) was touched, adding a pageSetup node.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
